# PTOA_ImageAnalysisPipeline.pptx — rename the M-file "mri_fitp.m" to
# "mri_fitps.m" everywhere it (and its derived file names) appear on the
# slide, per the commit:
#   "Update image procressing pipeline
#    Updated for current analysis M-file:  mri_fitps.m."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Replace the first occurrence of $OldText with $NewText inside a shape's
# text, using Characters() so only the targeted run(s) are touched and the
# rest of the shape's runs/formatting are left alone.
# (Positional parameters only -- this host's PowerShell subset does not
# bind named/splatted `-Param value` arguments reliably.)
function Replace-ShapeSubstring($Shape, $OldText, $NewText) {
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $pos = $full.IndexOf($OldText)
    if ($pos -lt 0) {
        return $false
    }
    $chars = $tr.Characters($pos + 1, $OldText.Length)
    $chars.Text = $NewText
    return $true
}

# Find the (first) shape on the slide whose text contains $Needle.
function Find-ShapeContaining($Slide, $Needle) {
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $candidate = $Slide.Shapes.Item($i)
        if ($candidate.HasTextFrame) {
            $t = $candidate.TextFrame.TextRange.Text
            if ($t.IndexOf($Needle) -ge 0) {
                return $candidate
            }
        }
    }
    return $null
}

# --- Slide text boxes referencing the M-file / its generated files ---

# 1) "10. mri_fitp.m" -> "10. mri_fitps.m"  (TextBox 16 - the M-file step)
$shapeM = Find-ShapeContaining $s "10. mri_fitp.m"
$ok1 = Replace-ShapeSubstring $shapeM "mri_fitp.m" "mri_fitps.m"

# 2) "mri_fitp*.ps/.pdf" -> "mri_fitps*.ps/.pdf"  (TextBox 18 - Plot Files)
$shapePlots = Find-ShapeContaining $s "mri_fitp*"
$ok2 = Replace-ShapeSubstring $shapePlots "mri_fitp" "mri_fitps"

# 3) "mri_fitp.xlsx" -> "mri_fitps.xlsx"  (TextBox 29 - xlsx Data File)
$shapeXlsx = Find-ShapeContaining $s "mri_fitp.xlsx"
$ok3 = Replace-ShapeSubstring $shapeXlsx "mri_fitp.xlsx" "mri_fitps.xlsx"

# 4) "mri_fitp.mat" -> "mri_fitps.mat"  (TextBox 35 - mat Data File)
#    The source author ended up with this run split into four pieces
#    (mri / _ / fitps. / mat); reproduce the same run boundaries by
#    editing "fitp." -> "fitps." first (splits into mri_ | fitps. | mat)
#    and then splitting "mri_" into "mri" | "_".
$shapeMat = Find-ShapeContaining $s "mri_fitp.mat"
$ok4a = Replace-ShapeSubstring $shapeMat "fitp." "fitps."
$ok4b = Replace-ShapeSubstring $shapeMat "_" "_"

Write-Host "Replacements applied (M/.ps.pdf/.xlsx/.mat):" $ok1 $ok2 $ok3 ($ok4a -and $ok4b)
Write-Host "M-file shape text:" $shapeM.TextFrame.TextRange.Text
Write-Host "Plot Files shape text:" $shapePlots.TextFrame.TextRange.Text
Write-Host "xlsx shape text:" $shapeXlsx.TextFrame.TextRange.Text
Write-Host "mat shape text:" $shapeMat.TextFrame.TextRange.Text
